$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.488.02"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'3.551.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'140.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("D7").Value = "'3.550.58"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Value = "'7.15"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("D13").Value = "'4.153.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("E14").Value = "  +4.97%  "
$ws.Range("D15").Value = "'27.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "'3.543.08"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "'65.322.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'10.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "'14.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.48%  "
$ws.Range("D22").Value = "'396.84"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "'0.571"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.27%  "
$ws.Range("D24").Value = "'74.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "'3.687.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +10.20%  "
$ws.Range("E28").Value = "  +8.24%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "'2.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'8.29"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").Value = "'3.564.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("D33").Value = "'24.04"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.11%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.147"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("D38").Value = "'168.43"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").Value = "'4.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").Value = "'0.0806"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("D42").Value = "'0.826"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").Value = "'26.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +20.67%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("E47").Value = "  +9.37%  "
$ws.Range("D48").Value = "'1.68"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").Value = "'6.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.19%  "
$ws.Range("D50").Value = "'2.389.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.40%  "
$ws.Range("D51").Value = "'2.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.38%  "
